$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before B for the invoice number field.
#        This shifts the old B..I (account..amount) one column right to C..J.
$ws.Columns.Item(2).Insert()

# --- 2. Populate the three new invoiceNumber placeholder cells.
$ws.Range("B1").Value = "{d.i18n.invoiceNumber}"
$ws.Range("B2").Value = "{d.objs[i].invoiceNumber}"
$ws.Range("B3").Value = "{d.objs[i+1].invoiceNumber}"

# --- 3. Fix the pre-existing "amount" placeholder typos (lowercase i -> I),
#        matching the already-corrected {d.objs[I].order}/{d.objs[I].status} cells.
#        These cells now live in column J after the column insert above.
$ws.Range("J2").Value = "{d.objs[I].amount}"
$ws.Range("J3").Value = "{d.objs[I+1].amount}"

# --- 4. Restore the per-column widths that the column insert carried along
#        with the shifted content, so each column position keeps (roughly)
#        the width it had before the edit.
$ws.Columns.Item(1).ColumnWidth = 37
$ws.Columns.Item(2).ColumnWidth = 15.83
$ws.Columns.Item(3).ColumnWidth = 35
$ws.Columns.Item(4).ColumnWidth = 9.67
$ws.Columns.Item(5).ColumnWidth = 21.17
$ws.Columns.Item(6).ColumnWidth = 28.5
$ws.Columns.Item(7).ColumnWidth = 25.17
$ws.Columns.Item(8).ColumnWidth = 13.67
$ws.Columns.Item(9).ColumnWidth = 13.67
$ws.Columns.Item(10).ColumnWidth = 13.67

# --- 5. The newly inserted column B picked up the left-edge border style from
#        column A in row 4 (the thin divider row under the headers). It should
#        instead look like the other interior columns of that row.
$ws.Range("C4").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 6. Flatten the per-column left/right border styling of the data rows
#        (5-10) into one uniform interior style across every column.
$ws.Range("C5").Copy()
$ws.Range("A5:J10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 7. Extend the sheet with five more blank, identically formatted data
#        rows (11-15), matching the rest of the data area.
$ws.Range("A5:J5").Copy()
$ws.Range("A11:J15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A11:J15").RowHeight = 15

$ws.Range("A1").Select()
